# Corrections Colm - 1
#
# 1) Update the cached "datetimeFigureOut" footer/date field text from
#    15/04/2019 -> 01/03/2020 on the slide master and every slide layout
#    (this is what PowerPoint does across Insert > Header & Footer > Apply
#    to All).
# 2) Fix the "Warwick" -> "WARP" typo in the "Warwick Model" shape label
#    on slide 1 (only the first run of that text box changes).

$p = $ppt.ActivePresentation

$oldDate = "15/04/2019"
$newDate = "01/03/2020"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# Slide 1: "Warwick Model" shape -> "WARP Model" (first run only).
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $firstRun = $shp.TextFrame.TextRange.Runs(1)
        if ($firstRun.Text -eq "Warwick") {
            $firstRun.Text = "WARP"
        }
    }
}
